# "qpsk 600 l26 revert to previous best"
#
# Revert the per-branch "scale factor" inputs on both loop-filter
# calculator sheets back to their previous-best values, and tidy up a
# stray leftover formatted row on the "Branch LPF" sheet. All of the
# dependent formulas (B20/B21 on sheet1, B22/B23 on sheet2) recalc
# automatically from these inputs.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "LoopFilter LPF"
$ws2 = $wb.Worksheets.Item(2)   # "Branch LPF"

# --- LoopFilter LPF: scale factor 4 -> 8 -------------------------------
$ws1.Range("B16").Value = 8

# --- Branch LPF: scale factor 4 -> 1 -----------------------------------
$ws2.Range("B18").Value = 1

# --- Branch LPF: drop the stray formatted row 26 (A26:B26) left behind
#     from earlier edits; it has no content, only formatting.
$ws2.Rows.Item(26).Delete()

# --- Restore the selections / active sheet seen in the saved file ------
$ws2.Range("B19").Select()
$ws1.Activate()
$ws1.Range("B17").Select()
